# Plancheck/protocol_check/prostate.xlsx
# "check the presence of structures + HU value + volume value"
#
# Adds HU / vol min / vol max columns (header row) to the "Clinical
# Structures" and "opt structures" sheets, fills in an example HU value
# (Vessie) and vol min/max values (Rectum) on "Clinical Structures", and
# renames/re-adds the Orfit couch-structure rows (capitalisation fix +
# new vol min / vol max header) on "couch_structures".

$wb = $excel.ActiveWorkbook

# --- couch_structures: fix Orfit capitalisation first so the shared------
# string table regenerates "zzInt_ORFIT" / "zzExt_ORFIT" in that order,
# freeing up the old "zzExt_Orfit" / "zzInt_Orfit" strings.
$ws4 = $wb.Worksheets.Item("couch_structures")
$ws4.Range("A5").Value = "zzInt_ORFIT"
$ws4.Range("A4").Value = "zzExt_ORFIT"

# --- Clinical Structures: add HU / vol min / vol max header + sample ---
# values (so "HU" becomes the next reused shared string, then the new
# "vol min" / "vol max" strings get appended).
$ws2 = $wb.Worksheets.Item("Clinical Structures")
$ws2.Range("B1").Value = "HU"
$ws2.Range("C1").Value = "vol min"
$ws2.Range("D1").Value = "vol max"
$ws2.Range("C17").Value = 30
$ws2.Range("D17").Value = 50
$ws2.Range("B23").Value = 0

# --- opt structures: add the same HU / vol min / vol max header -------
$ws3 = $wb.Worksheets.Item("opt structures")
$ws3.Range("B1").Value = "HU"
$ws3.Range("C1").Value = "vol min"
$ws3.Range("D1").Value = "vol max"

# --- couch_structures: add the same HU / vol min / vol max header -----
$ws4.Range("B1").Value = "HU"
$ws4.Range("C1").Value = "vol min"
$ws4.Range("D1").Value = "vol max"

# --- Selections / active sheet -----------------------------------------
# Final state: "Clinical Structures" is the active tab (activeTab=1),
# selection on C17; "opt structures" selection parked at F9 (no more
# topLeftCell scroll anchor); "couch_structures" selection on C1:D1 and
# no longer the active tab.
$ws3.Activate() | Out-Null
$ws3.Range("F9").Select() | Out-Null

$ws4.Activate() | Out-Null
$ws4.Range("C1:D1").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("C17").Select() | Out-Null
